# Apply the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Mon Oct 23 02:43:10 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text such as "30.285.82" or "0.0620" that Excel would
# otherwise auto-convert to a number (dropping separators / trailing zeros),
# so force the cells to Text format before writing the new values.
$dCells = @("D2","D3","D5","D8","D10","D11","D12","D13","D14","D15","D17","D18","D19","D22","D23","D24","D25","D26","D29","D31","D32","D34","D38","D39","D40","D41","D42","D43","D46","D47","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.285.82"
$ws.Range("E2").Value = "  +1.34%  "

$ws.Range("D3").Value = "1.675.58"
$ws.Range("E3").Value = "  +2.99%  "

$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").Value = "219.84"
$ws.Range("E5").Value = "  +2.53%  "

$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("D8").Value = "29.84"
$ws.Range("E8").Value = "  +0.93%  "

$ws.Range("E9").Value = "  +2.36%  "

$ws.Range("D10").Value = "0.0620"
$ws.Range("E10").Value = "  +1.17%  "

$ws.Range("D11").Value = "0.0907"
$ws.Range("E11").Value = "  -0.98%  "

$ws.Range("D12").Value = "1.914.98"
$ws.Range("E12").Value = "  +2.90%  "

$ws.Range("D13").Value = "10.77"
$ws.Range("E13").Value = "  +19.60%  "

$ws.Range("D14").Value = "0.621"
$ws.Range("E14").Value = "  +9.18%  "

$ws.Range("D15").Value = "1.649.13"
$ws.Range("E15").Value = "  +1.35%  "

$ws.Range("E16").Value = "  +3.17%  "

$ws.Range("D17").Value = "30.291.13"
$ws.Range("E17").Value = "  +1.19%  "

$ws.Range("D18").Value = "65.77"
$ws.Range("E18").Value = "  +1.51%  "

$ws.Range("D19").Value = "247.90"
$ws.Range("E19").Value = "  +0.75%  "

$ws.Range("E20").Value = "  +2.40%  "

$ws.Range("E21").Value = "  -0.21%  "

$ws.Range("D22").Value = "4.32"
$ws.Range("E22").Value = "  +4.84%  "

$ws.Range("D23").Value = "10.07"
$ws.Range("E23").Value = "  +5.15%  "

$ws.Range("D24").Value = "2.21"
$ws.Range("E24").Value = "  +4.41%  "

$ws.Range("D25").Value = "158.83"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("D26").Value = "15.90"
$ws.Range("E26").Value = "  +1.43%  "

$ws.Range("E27").Value = "  +0.49%  "

$ws.Range("E28").Value = "  +2.80%  "

$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("E30").Value = "  +2.22%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.14"
$ws.Range("E31").Value = "  +0.36%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "3.47"
$ws.Range("E32").Value = "  +3.85%  "

$ws.Range("E33").Value = "  +3.59%  "

$ws.Range("D34").Value = "1.481.33"
$ws.Range("E34").Value = "  +3.77%  "

$ws.Range("E35").Value = "  +5.82%  "

$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("E37").Value = "  +5.09%  "

$ws.Range("B38").Value = "Aave"
$ws.Range("C38").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D38").Value = "79.94"
$ws.Range("E38").Value = "  +15.56%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.591"
$ws.Range("E39").Value = "  +6.63%  "

$ws.Range("D40").Value = "2.69"
$ws.Range("E40").Value = "  -6.58%  "

$ws.Range("D41").Value = "2.29"
$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("D42").Value = "0.860"
$ws.Range("E42").Value = "  +3.48%  "

$ws.Range("D43").Value = "2.02"
$ws.Range("E43").Value = "  +3.03%  "

$ws.Range("E44").Value = "  +1.01%  "

$ws.Range("E45").Value = "  -3.18%  "

$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("D47").Value = "52.23"
$ws.Range("E47").Value = "  -5.01%  "

$ws.Range("E48").Value = "  +0.22%  "

$ws.Range("D49").Value = "1.808.83"
$ws.Range("E49").Value = "  +2.25%  "

$ws.Range("D50").Value = "95.37"
$ws.Range("E50").Value = "  +6.40%  "

$ws.Range("D51").Value = "0.0₆0117"
$ws.Range("E51").Value = "  +9.53%  "

Write-Host "Applied cryptos.xlsx price/volume refresh"
